$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.704.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.794.35"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.34%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.793.50"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.35%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("E10").Value = "  +0.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.31"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.61%  "
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("E13").Value = "  -1.63%  "
$ws.Range("E14").Value = "  +0.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.428.73"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.802.03"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.652.52"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.37%  "
$ws.Range("E18").Value = "  +3.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.04"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.79%  "
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "459.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("E24").Value = "  +5.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.29"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.06"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.51%  "
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("E29").Value = "  -0.43%  "
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("E31").Value = "  +4.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.66"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.09"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.734.88"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.100"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.38"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.84%  "
$ws.Range("E39").Value = "  +0.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.991"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.26%  "
$ws.Range("E41").Value = "  +0.34%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "44.18"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "48.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.78%  "
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "149.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.13%  "
$ws.Range("E48").Value = "  -1.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "393.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "26.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.94%  "
$ws.Range("E51").Value = "  -4.03%  "
